$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.2182345178183951
$ws.Range("B1").Value = 0.21804058587594
$ws.Range("A2").Value = -0.15727086447088112
$ws.Range("B2").Value = 0.15685396304133015
$ws.Range("A3").Value = -0.1071485688817031
$ws.Range("B3").Value = 0.10697114947461372
$ws.Range("A4").Value = -0.09897114953202113
$ws.Range("B4").Value = 0.09858128091380891
$ws.Range("A5").Value = -0.09558128094641916
$ws.Range("B5").Value = 0.09426242480678493
$ws.Range("A6").Value = -0.04033661058318749
$ws.Range("B6").Value = 0.04001328858917752
$ws.Range("A7").Value = -0.030013288672153582
$ws.Range("B7").Value = 0.02994449497593754
$ws.Range("A8").Value = -0.019944495061059886
$ws.Range("B8").Value = 0.019848439018641084
$ws.Range("A9").Value = -0.0178484390612903
$ws.Range("B9").Value = 0.017777640625388003
$ws.Range("A10").Value = -0.01577764067030074
$ws.Range("B10").Value = 0.01577469520703545
$ws.Range("A11").Value = -0.012774695257741442
$ws.Range("B11").Value = 0.012769695076466192
$ws.Range("A12").Value = -0.009269695130392552
$ws.Range("B12").Value = 0.009240882387173688
$ws.Range("A13").Value = -0.005740882442913708
$ws.Range("B13").Value = 0.005733985863359159
$ws.Range("A14").Value = 0.0022660140551131747
$ws.Range("B14").Value = -0.002266056368161351
$ws.Range("A15").Value = 0.003266056325972322
$ws.Range("B15").Value = -0.0032679850042036662
$ws.Range("A16").Value = -0.006033963824158217
$ws.Range("B16").Value = 0.006003423258485352
$ws.Range("A17").Value = -0.004003423307185727
$ws.Range("B17").Value = 0.003999999939794385
$ws.Range("A18").Value = -0.0547697214402767
$ws.Range("B18").Value = 0.054666415890579856
$ws.Range("A19").Value = -0.05066641591524457
$ws.Range("B19").Value = 0.049910654173991986
$ws.Range("A20").Value = -0.045910654207485635
$ws.Range("B20").Value = 0.045697554753706626
$ws.Range("A21").Value = -0.004005798857315668
$ws.Range("B21").Value = 0.003999999964524825
$ws.Range("A22").Value = -0.04570539419125019
$ws.Range("B22").Value = 0.04549411297671213
$ws.Range("A23").Value = -0.04049411301580186
$ws.Range("B23").Value = 0.040098037145776644
$ws.Range("A24").Value = -0.020098037272974878
$ws.Range("B24").Value = 0.019999999871133767
$ws.Range("A25").Value = -0.05192581426039311
$ws.Range("B25").Value = 0.05188920117455531
$ws.Range("A26").Value = -0.04938920121423607
$ws.Range("B26").Value = 0.049345350170680646
$ws.Range("A27").Value = -0.04684535021154179
$ws.Range("B27").Value = 0.04660151496849041
$ws.Range("A28").Value = -0.04460151501197451
$ws.Range("B28").Value = 0.04444901014741287
$ws.Range("A29").Value = -0.03744901022286928
$ws.Range("B29").Value = 0.0374156074441272
$ws.Range("A30").Value = 0.022584392180328816
$ws.Range("B30").Value = -0.02261810925610641
$ws.Range("A31").Value = 0.029618109180864494
$ws.Range("B31").Value = -0.029641039974835692
$ws.Range("A32").Value = -0.004001239050550609
$ws.Range("B32").Value = 0.003999999942379873
